# "OC 4 - beregnTvaertkraft er lavet"
#
# Nada logs more time-tracking entries in the "Tidsregistrering" sheet:
#   - a couple of earlier entries had their duration text corrected
#     (missing " min." suffix, and a "FURSP" -> "FURPS" typo fix)
#   - row 16 (previously an empty placeholder row) gets filled in with
#     a new logged activity ("Skabelon til OC")
#   - two brand-new rows are inserted for two more activities
#     ("OC1 for angivVaegt" and "OC3 for beregnNormalkraft")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Make room for the two new log entries: push rows 17+ down by two rows,
# copying the formatting (styles/number formats) from row 16 above.
$ws.Rows("17:18").Insert() | Out-Null

# --- Correct a few existing entries' duration text / spelling ---
$ws.Range("I11").Value = "1 time : 00 min."
$ws.Range("I12").Value = "2 time : 50 min."
$ws.Range("F13").Value = "FURPS"
$ws.Range("I13").Value = "0 time : 30 min."
$ws.Range("I14").Value = "0 time : 30 min."
$ws.Range("I15").Value = "0 time : 20 min."

# --- Row 16: fill in the previously-blank activity row ---
$ws.Range("E16").Value = "Requirements Specifier"
$ws.Range("F16").Value = "Skabelon til OC"
$ws.Range("G16").Value = 0.54166666666666663
$ws.Range("H16").Value = 0.54866898148148147
$ws.Range("I16").Value = "0 time : 05 min."

# --- Row 17: new activity entry ---
$ws.Range("A17").Value = 42802
$ws.Range("B17").Value = "NO"
$ws.Range("C17").Value = "Nada H. A. Omer"
$ws.Range("E17").Value = "Requirements Specifier"
$ws.Range("F17").Value = "OC1 for angivVaegt"
$ws.Range("G17").Value = 0.55208333333333337
$ws.Range("H17").Value = 0.57291666666666663
$ws.Range("I17").Value = "0 time : 30 min."

# --- Row 18: new activity entry ---
$ws.Range("A18").Value = 42802
$ws.Range("B18").Value = "NO"
$ws.Range("C18").Value = "Nada H. A. Omer"
$ws.Range("E18").Value = "Requirements Specifier"
$ws.Range("F18").Value = "OC3 for beregnNormalkraft"
$ws.Range("G18").Value = 0.57638888888888895
$ws.Range("H18").Value = 0.62847222222222221
$ws.Range("I18").Value = "1 time : 15 min."

# Leave the selection where the author left off editing.
$ws.Range("B23").Select() | Out-Null
